$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 502103.9
$ws.Range("J6").Value = 4200
$ws.Range("L6").Value = 12600
$ws.Range("N6").Value = -12824
$ws.Range("H12").Value = 125
$ws.Range("I12").Value = 125
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 125
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = 45
$ws.Range("N12").Value = $null
$ws.Range("H21").Value = 6812.6665
$ws.Range("I21").Value = 400
$ws.Range("K21").Value = 400
$ws.Range("M21").Value = 68
$ws.Range("H23").Value = 6812.6665
$ws.Range("I23").Value = 400
$ws.Range("K23").Value = 400
$ws.Range("M23").Value = -166
$ws.Range("H43").Value = 2100
$ws.Range("I43").Value = 2100
$ws.Range("K43").Value = 2100
$ws.Range("M43").Value = -2031
$ws.Range("H62").Value = 1999.5
$ws.Range("I62").Value = 2000
$ws.Range("J62").Value = 1999
$ws.Range("K62").Value = 2000
$ws.Range("L62").Value = 1999
$ws.Range("M62").Value = -1376
$ws.Range("N62").Value = -3247
$ws.Range("H65").Value = 1999.5
$ws.Range("I65").Value = 2000
$ws.Range("J65").Value = 1999
$ws.Range("K65").Value = 10000
$ws.Range("L65").Value = 9995
$ws.Range("M65").Value = -6880
$ws.Range("N65").Value = -16235
$ws.Range("H116").Value = 7660
$ws.Range("I116").Value = 6766.6665
$ws.Range("K116").Value = 6766.6665
$ws.Range("M116").Value = -3324.6665
$ws.Range("H125").Value = 3929.9
$ws.Range("I125").Value = 1649.5
$ws.Range("J125").Value = 4500
$ws.Range("K125").Value = 14845.5
$ws.Range("L125").Value = 40500
$ws.Range("M125").Value = -12385.5
$ws.Range("N125").Value = -45420
$ws.Range("H127").Value = 2365.6667
$ws.Range("I127").Value = 1048.5
$ws.Range("J127").Value = 5000
$ws.Range("K127").Value = 3145.5
$ws.Range("L127").Value = 15000
$ws.Range("N127").Value = -24920
$ws.Range("M127").Value = 1814.5
$ws.Range("H131").Value = 2663.1667
$ws.Range("I131").Value = 1195.8
$ws.Range("J131").Value = 10000
$ws.Range("K131").Value = 3587.4
$ws.Range("L131").Value = 30000
$ws.Range("M131").Value = 1452.6
$ws.Range("N131").Value = -40080

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 262.5
$ws.Range("I3").Value = 262.5
$ws.Range("K3").Value = 262.5
$ws.Range("M3").Value = -147.5
$ws.Range("H22").Value = 2000
$ws.Range("I22").Value = 2000
$ws.Range("K22").Value = 2000
$ws.Range("M22").Value = -1701
$ws.Range("H35").Value = 1750
$ws.Range("I35").Value = 1000
$ws.Range("J35").Value = 2500
$ws.Range("K35").Value = 1000
$ws.Range("L35").Value = 2500
$ws.Range("N35").Value = -3312
$ws.Range("M35").Value = -594
$ws.Range("H113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").Value = $null
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").Value = $null

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 820.2
$ws.Range("I7").Value = 1215.6666
$ws.Range("J7").Value = 227
$ws.Range("K7").Value = 1215.6666
$ws.Range("L7").Value = 227
$ws.Range("M7").Value = -1102.6666
$ws.Range("N7").Value = -453
$ws.Range("H22").Value = 3433.6667
$ws.Range("I22").Value = 3433.6667
$ws.Range("K22").Value = 3433.6667
$ws.Range("M22").Value = -3260.6667
$ws.Range("H99").Value = 3900
$ws.Range("I99").Value = 4250
$ws.Range("K99").Value = 4250
$ws.Range("M99").Value = -2752

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 373.85715
$ws.Range("I2").Value = 269.5
$ws.Range("K2").Value = 269.5
$ws.Range("M2").Value = -156.5
$ws.Range("H16").Value = 1749.5
$ws.Range("I16").Value = 1749.5
$ws.Range("K16").Value = 1749.5
$ws.Range("M16").Value = -1462.5
$ws.Range("H31").Value = 11390
$ws.Range("I31").Value = 4012
$ws.Range("K31").Value = 4012
$ws.Range("M31").Value = -3717
$ws.Range("H34").Value = 11390
$ws.Range("I34").Value = 4012
$ws.Range("K34").Value = 4012
$ws.Range("M34").Value = -3810
$ws.Range("H35").Value = 275
$ws.Range("I35").Value = 275
$ws.Range("K35").Value = 275
$ws.Range("M35").Value = 19
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("M50").Value = $null
$ws.Range("H62").Value = 5500
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 5500
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 5500
$ws.Range("N62").Value = -6748
$ws.Range("M62").Value = $null
$ws.Range("H65").Value = 5500
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 5500
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 27500
$ws.Range("N65").Value = -33740
$ws.Range("M65").Value = $null
$ws.Range("H113").Value = 1749.5
$ws.Range("I113").Value = 1749.5
$ws.Range("K113").Value = 1749.5
$ws.Range("M113").Value = 420.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 213.83333
$ws.Range("J15").Value = 213.83333
$ws.Range("L15").Value = 641.49999
$ws.Range("N15").Value = -921.49999
$ws.Range("H16").Value = 11749.875
$ws.Range("J16").Value = 11749.875
$ws.Range("L16").Value = 35249.625
$ws.Range("N16").Value = -35595.625
$ws.Range("H138").Value = 2299.6667
$ws.Range("I138").Value = 2299.6667
$ws.Range("K138").Value = 6899.000100000001
$ws.Range("M138").Value = -1759.000100000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 61.18182
$ws.Range("I2").Value = 68
$ws.Range("J2").Value = 49.25
$ws.Range("K2").Value = 68
$ws.Range("L2").Value = 49.25
$ws.Range("M2").Value = 45
$ws.Range("N2").Value = -275.25
$ws.Range("H126").Value = 1468.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 461
$ws.Range("I9").Value = 339.5
$ws.Range("K9").Value = 339.5
$ws.Range("M9").Value = -115.5
$ws.Range("H13").Value = 2577
$ws.Range("I13").Value = 618.2857
$ws.Range("K13").Value = 618.2857
$ws.Range("M13").Value = -478.2857
$ws.Range("H14").Value = 5001
$ws.Range("J14").Value = 5001
$ws.Range("L14").Value = 5001
$ws.Range("N14").Value = -5345
$ws.Range("H16").Value = 4950
$ws.Range("J16").Value = 4900
$ws.Range("L16").Value = 4900
$ws.Range("N16").Value = -5240
$ws.Range("H31").Value = 800
$ws.Range("J31").Value = 800
$ws.Range("L31").Value = 800
$ws.Range("N31").Value = -1296
$ws.Range("H82").Value = 2389
$ws.Range("J82").Value = 2333.6667
$ws.Range("L82").Value = 2333.6667
$ws.Range("N82").Value = -3055.6667
$ws.Range("H85").Value = 2389
$ws.Range("J85").Value = 2333.6667
$ws.Range("L85").Value = 2333.6667
$ws.Range("N85").Value = -4829.6667
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").Value = $null

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H19").Value = 12000
$ws.Range("J19").Value = 12000
$ws.Range("L19").Value = 12000
$ws.Range("N19").Value = -12348
$ws.Range("H23").Value = 353
$ws.Range("I23").Value = 353
$ws.Range("K23").Value = 353
$ws.Range("M23").Value = -124
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("M24").Value = $null
$ws.Range("H33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").Value = $null
$ws.Range("H36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").Value = $null
$ws.Range("H136").Value = 13001.25
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").Value = $null
